$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 74.19408395296784
$ws.Cells.Item(2, 3).Value = 235.6160305053561
$ws.Cells.Item(3, 2).Value = 89.43834202308427
$ws.Cells.Item(3, 3).Value = 322.8997759563496
$ws.Cells.Item(4, 2).Value = 81.40379101139183
$ws.Cells.Item(4, 3).Value = 295.5660860094558
$ws.Cells.Item(5, 2).Value = 86.83916451945294
$ws.Cells.Item(5, 3).Value = 321.1505416764683
$ws.Cells.Item(6, 2).Value = 12.89914685508182
$ws.Cells.Item(6, 3).Value = 48.79759880703651
$ws.Cells.Item(10, 2).Value = 33.88210938269473
$ws.Cells.Item(10, 3).Value = 144.5162077920627
$ws.Cells.Item(11, 2).Value = 70.47520533542121
$ws.Cells.Item(11, 3).Value = 273.6202019474335
$ws.Cells.Item(12, 2).Value = 76.36126498340384
$ws.Cells.Item(12, 3).Value = 326.9761422393509
$ws.Cells.Item(13, 2).Value = 67.29923937340671
$ws.Cells.Item(13, 3).Value = 320.7905813088454
$ws.Cells.Item(14, 2).Value = 68.12926509896874
$ws.Cells.Item(14, 3).Value = 284.4968567438901
$ws.Cells.Item(15, 2).Value = 76.68185093707488
$ws.Cells.Item(15, 3).Value = 295.1250822512637
$ws.Cells.Item(16, 2).Value = 71.72916689551374
$ws.Cells.Item(16, 3).Value = 307.6276947554144
$ws.Cells.Item(17, 2).Value = 49.52876982530729
$ws.Cells.Item(17, 3).Value = 207.7216895924848
$ws.Cells.Item(18, 2).Value = 62.3870836582288
$ws.Cells.Item(18, 3).Value = 282.1812468481443
$ws.Cells.Item(19, 2).Value = 17.3163181737962
$ws.Cells.Item(19, 3).Value = 85.2829193704211
$ws.Cells.Item(332, 2).Value = 49.53641628222108
$ws.Cells.Item(332, 3).Value = 195.6177744368306
$ws.Cells.Item(333, 2).Value = 67.65603271884665
$ws.Cells.Item(333, 3).Value = 265.4257390016767
$ws.Cells.Item(334, 2).Value = 53.43321041568568
$ws.Cells.Item(334, 3).Value = 180.5277599601374
$ws.Cells.Item(335, 2).Value = 51.88458076800471
$ws.Cells.Item(335, 3).Value = 185.3590561792269
$ws.Cells.Item(336, 2).Value = 22.99160156496697
$ws.Cells.Item(336, 3).Value = 82.42496309863617
$ws.Cells.Item(337, 2).Value = 90.35570131241732
$ws.Cells.Item(337, 3).Value = 312.9842346499702
$ws.Cells.Item(338, 2).Value = 66.47221519592834
$ws.Cells.Item(338, 3).Value = 236.0801356833313
$ws.Cells.Item(339, 2).Value = 63.49660768735646
$ws.Cells.Item(339, 3).Value = 227.3674781209431
$ws.Cells.Item(340, 2).Value = 71.00953652726186
$ws.Cells.Item(340, 3).Value = 292.6441106717706
$ws.Cells.Item(393, 2).Value = 52.3238815982878
$ws.Cells.Item(393, 3).Value = 184.4766616097709
$ws.Cells.Item(394, 2).Value = 81.74441741544337
$ws.Cells.Item(394, 3).Value = 284.273573943435
$ws.Cells.Item(395, 2).Value = 66.50176335058171
$ws.Cells.Item(395, 3).Value = 255.997714289286
$ws.Cells.Item(396, 2).Value = 54.39357174783888
$ws.Cells.Item(396, 3).Value = 208.9845211338032
$ws.Cells.Item(397, 2).Value = 55.33902895397608
$ws.Cells.Item(397, 3).Value = 211.6259761430322
$ws.Cells.Item(415, 2).Value = 89.72959203905191
$ws.Cells.Item(415, 3).Value = 270.8278930788219
$ws.Cells.Item(416, 2).Value = 89.95104541581306
$ws.Cells.Item(416, 3).Value = 251.9425414092616
$ws.Cells.Item(417, 2).Value = 95.8267613971673
$ws.Cells.Item(417, 3).Value = 319.263771972892
$ws.Cells.Item(418, 2).Value = 112.96153866474
$ws.Cells.Item(418, 3).Value = 341.0152590276041
$ws.Cells.Item(419, 2).Value = 57.09796487215972
$ws.Cells.Item(419, 3).Value = 186.1405822876012
$ws.Cells.Item(454, 2).Value = 41.09573104696092
$ws.Cells.Item(454, 3).Value = 159.770758626355
$ws.Cells.Item(455, 2).Value = 80.68603951886227
$ws.Cells.Item(455, 3).Value = 310.0930647615996
$ws.Cells.Item(456, 2).Value = 73.16593037566265
$ws.Cells.Item(456, 3).Value = 288.1524420535704
$ws.Cells.Item(457, 2).Value = 40.61371398592992
$ws.Cells.Item(457, 3).Value = 152.3069586493511
$ws.Cells.Item(458, 2).Value = 43.84693099188979
$ws.Cells.Item(458, 3).Value = 171.1594880868015
$ws.Cells.Item(510, 2).Value = 23.46702885390863
$ws.Cells.Item(510, 3).Value = 94.03088363335192
$ws.Cells.Item(511, 2).Value = 71.49740029740083
$ws.Cells.Item(511, 3).Value = 301.6934536168001
$ws.Cells.Item(512, 2).Value = 81.13655678641017
$ws.Cells.Item(512, 3).Value = 307.490211741934
$ws.Cells.Item(513, 2).Value = 86.95960832714304
$ws.Cells.Item(513, 3).Value = 309.3345309898328
$ws.Cells.Item(514, 2).Value = 87.81513156460471
$ws.Cells.Item(514, 3).Value = 299.7865204695406
$ws.Cells.Item(515, 2).Value = 90.34599235316769
$ws.Cells.Item(515, 3).Value = 314.4756660099475
$ws.Cells.Item(516, 2).Value = 47.05255936190578
$ws.Cells.Item(516, 3).Value = 149.3019578533691
$ws.Cells.Item(517, 2).Value = 77.77085003318966
$ws.Cells.Item(517, 3).Value = 251.7955005061134
$ws.Cells.Item(518, 2).Value = 72.25796910256551
$ws.Cells.Item(518, 3).Value = 232.5003231081802
$ws.Cells.Item(519, 2).Value = 19.28184452671116
$ws.Cells.Item(519, 3).Value = 62.06945016401679
$ws.Cells.Item(580, 2).Value = 65.22965875316933
$ws.Cells.Item(580, 3).Value = 264.1628533872337
$ws.Cells.Item(581, 2).Value = 78.63398330660458
$ws.Cells.Item(581, 3).Value = 311.2252101021177
$ws.Cells.Item(582, 2).Value = 66.90199116984282
$ws.Cells.Item(582, 3).Value = 272.1300033031059
$ws.Cells.Item(583, 2).Value = 71.65319751390444
$ws.Cells.Item(583, 3).Value = 313.1579700885604
$ws.Cells.Item(584, 2).Value = 43.94523252154485
$ws.Cells.Item(584, 3).Value = 159.7777856250967
$ws.Cells.Item(606, 2).Value = 35.94738360435686
$ws.Cells.Item(606, 3).Value = 131.2784698312382
$ws.Cells.Item(607, 2).Value = 68.58709693293379
$ws.Cells.Item(607, 3).Value = 259.9260202829755
$ws.Cells.Item(608, 2).Value = 58.7103703390521
$ws.Cells.Item(608, 3).Value = 243.0631912795414
$ws.Cells.Item(609, 2).Value = 49.90356429910386
$ws.Cells.Item(609, 3).Value = 195.2176278500575
$ws.Cells.Item(610, 2).Value = 43.8055337938171
$ws.Cells.Item(610, 3).Value = 177.4631878721851
